$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.933.19'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.559.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.03'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0856'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.782.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.560.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.70%  '
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.942.57'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  -1.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.416.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.62%  '
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  +0.29%  '
$ws.Range("E39").Value = '  +1.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("E44").Value = '  +2.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.65'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.696.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("E51").Value = '  -0.70%  '
